$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# --- Cells changing from a number to the "0" text placeholder (shared string) ---
# Donor cell C18 already holds text "0" with the correct style (13)
$ws.Range("C18").Copy($ws.Range("C16"))
$ws.Range("C18").Copy($ws.Range("G31"))

# --- Cells changing from a number to the "***.*" text placeholder (shared string) ---
# Donor cell L23 already holds text "***.*" with the correct style (13)
$ws.Range("L23").Copy($ws.Range("H31"))

# --- Cells changing from the "0"/"***.*" text placeholder to a real integer (style 14) ---
# Donor cell D18 already holds an integer with the correct style (14)
$ws.Range("D18").Copy($ws.Range("C25"))
$ws.Range("C25").Value = 2
$ws.Range("D18").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("D18").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("D18").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 2
$ws.Range("D18").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 2
$ws.Range("D18").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("D18").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1

# --- Cells changing from the "0"/"***.*" text placeholder to a real decimal (style 15) ---
# Donor cell E18 already holds a decimal with the correct style (15)
$ws.Range("E18").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("E18").Copy($ws.Range("H29"))
$ws.Range("H29").Value = -100
$ws.Range("E18").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("E18").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100

# --- Plain value updates (type/style unchanged) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 6
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = -27.027027027027
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 116
$ws.Range("J17").Value = 107
$ws.Range("K17").Value = 8.411214953271
$ws.Range("L17").Value = 1.754385964912
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = -66.666666666666
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -41.666666666666
$ws.Range("L18").Value = -38.235294117647
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 12
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = -7.692307692307
$ws.Range("I19").Value = 96
$ws.Range("J19").Value = 123
$ws.Range("K19").Value = -21.951219512195
$ws.Range("L19").Value = -22.580645161290
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 84
$ws.Range("J20").Value = 73
$ws.Range("K20").Value = 15.068493150684
$ws.Range("L20").Value = -2.325581395348
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 51
$ws.Range("G21").Value = 60
$ws.Range("H21").Value = -15
$ws.Range("I21").Value = 360
$ws.Range("J21").Value = 381
$ws.Range("K21").Value = -5.511811023622
$ws.Range("L21").Value = -10.447761194029
$ws.Range("C24").Value = 10
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 31
$ws.Range("G24").Value = 34
$ws.Range("H24").Value = -8.823529411764
$ws.Range("I24").Value = 220
$ws.Range("J24").Value = 267
$ws.Range("K24").Value = -17.602996254681
$ws.Range("L24").Value = -6.382978723404
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -42.857142857142
$ws.Range("I25").Value = 53
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = -18.461538461538
$ws.Range("L25").Value = 15.217391304347
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -36.363636363636
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 18.181818181818
$ws.Range("I26").Value = 194
$ws.Range("J26").Value = 188
$ws.Range("K26").Value = 3.191489361702
$ws.Range("L26").Value = 28.476821192053
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 133.333333333333
$ws.Range("L27").Value = 90.909090909090
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -22.222222222222
$ws.Range("L28").Value = -50
$ws.Range("J29").Value = 6
$ws.Range("J30").Value = 4

# --- Finally set the values for the cells that changed type to string placeholders ---
